# Apply updated cryptocurrency market data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as '1.002' or '328.50'. Excel's COM layer
# auto-converts plain-looking numeric text to actual numbers (dropping trailing
# zeros), so we force the Price column to Text format first to preserve the
# exact original string formatting.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Coin name swap (rows 41/42) ---
$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("B42").Value = 'TrustWalletToken'

# --- Link swap (rows 41/42) ---
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

# --- Price column (D) updates ---
$ws.Range("D2").Value = '27.770.79'
$ws.Range("D3").Value = '1.759.49'
$ws.Range("D4").Value = '1.002'
$ws.Range("D5").Value = '328.50'
$ws.Range("D7").Value = '0.4837'
$ws.Range("D8").Value = '0.3559'
$ws.Range("D9").Value = '43.17'
$ws.Range("D10").Value = '0.07519'
$ws.Range("D11").Value = '1.082'
$ws.Range("D13").Value = '20.55'
$ws.Range("D14").Value = '6.103'
$ws.Range("D15").Value = '7.119'
$ws.Range("D16").Value = '1.760.94'
$ws.Range("D17").Value = '93.28'
$ws.Range("D18").Value = '0.00001082'
$ws.Range("D19").Value = '0.06433'
$ws.Range("D20").Value = '1.001'
$ws.Range("D21").Value = '16.81'
$ws.Range("D22").Value = '5.798'
$ws.Range("D23").Value = '27.813.19'
$ws.Range("D24").Value = '11.10'
$ws.Range("D25").Value = '2.165'
$ws.Range("D26").Value = '163.58'
$ws.Range("D27").Value = '20.15'
$ws.Range("D28").Value = '1.960.84'
$ws.Range("D29").Value = '2.204'
$ws.Range("D30").Value = '122.64'
$ws.Range("D31").Value = '1.056'
$ws.Range("D32").Value = '0.09437'
$ws.Range("D33").Value = '3.655'
$ws.Range("D34").Value = '5.548'
$ws.Range("D35").Value = '0.02266'
$ws.Range("D36").Value = '11.61'
$ws.Range("D37").Value = '0.05976'
$ws.Range("D38").Value = '0.2064'
$ws.Range("D39").Value = '4.885'
$ws.Range("D40").Value = '0.6143'
$ws.Range("D41").Value = '1.434'
$ws.Range("D42").Value = '1.180'
$ws.Range("D43").Value = '7.749'
$ws.Range("D44").Value = '13.08'
$ws.Range("D45").Value = '3.732'
$ws.Range("D46").Value = '0.5791'
$ws.Range("D47").Value = '123.10'
$ws.Range("D48").Value = '1.926'
$ws.Range("D49").Value = '1.149'
$ws.Range("D50").Value = '0.06787'
$ws.Range("D51").Value = '71.74'

# --- Volume(1h) column (E) updates ---
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("E7").Value = '  +6.91%  '
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("E11").Value = '  -2.85%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("E15").Value = '  -2.12%  '
$ws.Range("E16").Value = '  -1.33%  '
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("E21").Value = '  -2.60%  '
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("E23").Value = '  -1.15%  '
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("E25").Value = '  +3.38%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("E29").Value = '  -3.15%  '
$ws.Range("E30").Value = '  -3.51%  '
$ws.Range("E31").Value = '  -5.26%  '
$ws.Range("E32").Value = '  +2.56%  '
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("E35").Value = '  -2.18%  '
$ws.Range("E36").Value = '  -3.25%  '
$ws.Range("E37").Value = '  -3.07%  '
$ws.Range("E38").Value = '  -2.04%  '
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("E40").Value = '  -3.83%  '
$ws.Range("E41").Value = '  +2.56%  '
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("E43").Value = '  -2.74%  '
$ws.Range("E44").Value = '  -2.52%  '
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("E46").Value = '  -2.70%  '
$ws.Range("E47").Value = '  -0.42%  '
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  -2.62%  '
$ws.Range("E51").Value = '  -1.96%  '
